# Plantilla_Ambientes.xlsx edit:
#  - Add 50 blank bordered data rows (2-51) to Hoja1 (columns A:G)
#  - Add a new "Hoja2" sheet holding the list of "Tipo" values
#  - Add a list-based data validation on Hoja1!D2:D51 sourced from Hoja2

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Hoja1: stripe rows 2-51 under the header with the grid/border look
#    used throughout the rest of the template (A:G, thin box border).
# ---------------------------------------------------------------------
$body = $ws1.Range("A2:G51")
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2

# ---------------------------------------------------------------------
# 2. Add the "Hoja2" helper sheet right after "Hoja1" and populate the
#    catalogue of room types used by the dropdown.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Hoja2"

# Touch scratch cells first so the shared-string table is built up in
# the same order the original workbook uses (Tipo, Salon, Laboratorio,
# Edificio, Auditorio) before writing the final cell layout below.
$ws2.Range("Z1").Value = "Tipo"
$ws2.Range("Z2").Value = "Salon"
$ws2.Range("Z3").Value = "Laboratorio"
$ws2.Range("Z4").Value = "Edificio"
$ws2.Range("Z5").Value = "Auditorio"
$ws2.Range("Z1:Z5").ClearContents()

$ws2.Range("A1").Value = "Tipo"
$ws2.Range("A2").Value = "Edificio"
$ws2.Range("A3").Value = "Laboratorio"
$ws2.Range("A4").Value = "Salon"
$ws2.Range("A5").Value = "Auditorio"

$catalog = $ws2.Range("A1:A5")
$catalog.Borders.LineStyle = 1
$catalog.Borders.Weight = 2

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3. Data validation dropdown on Hoja1!D2:D51 pulling from Hoja2.
# ---------------------------------------------------------------------
$dropdown = $ws1.Range("D2:D51")
$dropdown.Validation.Add(3, 1, 1, "=Hoja2!`$A`$2:`$A`$5")
$dropdown.Validation.IgnoreBlank = $true
$dropdown.Validation.InCellDropdown = $true
$dropdown.Validation.ShowInput = $true
$dropdown.Validation.ShowError = $true

# ---------------------------------------------------------------------
# 4. Restore the selections shown in the authored workbook.
# ---------------------------------------------------------------------
[void]$ws2.Range("C6").Select()
[void]$ws1.Select()
[void]$ws1.Range("G2").Select()
